# Generate Report for Handback
# Adds a second handed-back file (97f69da1-c9e3-46fd-943a-dd10195ca5ab.md) as a new
# row to each of the three report sheets (Overview, zh-cn, de-de), and refreshes the
# first file's (7c984b00... -> 6fcdd4e4...) generated-artifact names/timestamps.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$oldFile = "7c984b00-1d7f-429b-8bb7-65de61b23310"
$newFile1 = "6fcdd4e4-b5d7-4f1a-bbba-e5fa5aad9fc6"
$newFile2 = "97f69da1-c9e3-46fd-943a-dd10195ca5ab"

# ---------------------------------------------------------------------------
# 1. Grow each table by one row (auto-expands dimension, table ref & autofilter)
# ---------------------------------------------------------------------------
$ws1.ListObjects.Item(1).ListRows.Add() | Out-Null
$ws2.ListObjects.Item(1).ListRows.Add() | Out-Null
$ws3.ListObjects.Item(1).ListRows.Add() | Out-Null

# ---------------------------------------------------------------------------
# 2. Overview sheet
# ---------------------------------------------------------------------------
# Row 2 content refreshed (same file, new handback run)
$ws1.Range("A2").Value = "$newFile1.md"
$ws1.Range("C2").Value = ".md"
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws1.Range("G2").Value = "2016-08-23 19:05:26"

# Row 3 content (new file)
$ws1.Range("A3").Value = "$newFile2.md"
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-08-23 19:05:26"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks: rebuild in final order (B2, B3)
$ws1.Range("B2").Hyperlinks.Delete() | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "e2e\$newFile1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "e2e\$newFile2.md") | Out-Null

# ---------------------------------------------------------------------------
# 3. zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "$newFile1.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "e2e"
$ws2.Range("E2").Value = "ht"
$ws2.Range("F2").Value = "False"
$ws2.Range("G2").Value = "$newFile1.075c9c4432c64c1ca796c3749c5def87881cb935.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-23 19:05:21"
$ws2.Range("I2").Value = "$newFile1.md"
$ws2.Range("J2").Value = "$newFile1.075c9c4432c64c1ca796c3749c5def87881cb935.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-23 19:05:39"
$ws2.Range("L2").Value = ""
$ws2.Range("M2").Value = "True"
$ws2.Range("N2").Value = ""
$ws2.Range("O2").Value = "False"
$ws2.Range("P2").Value = ""

$ws2.Range("A3").Value = "$newFile2.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "$newFile2.419863fc8d9492e3b6dc1f1704fb5fb5fff97a10.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-23 19:05:21"
$ws2.Range("I3").Value = "$newFile2.md"
$ws2.Range("J3").Value = "$newFile2.419863fc8d9492e3b6dc1f1704fb5fb5fff97a10.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-23 19:05:39"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks: rebuild in final order (A2, I2, A3, I3)
$ws2.Range("A2").Hyperlinks.Delete() | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ee5cddfeb5f45a8f6b8515f07e33ab093ad473fa/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ee5cddfeb5f45a8f6b8515f07e33ab093ad473fa/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md") | Out-Null

# ---------------------------------------------------------------------------
# 4. de-de sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "$newFile1.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "e2e"
$ws3.Range("E2").Value = "ht"
$ws3.Range("F2").Value = "False"
$ws3.Range("G2").Value = "$newFile1.075c9c4432c64c1ca796c3749c5def87881cb935.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-23 19:05:26"
$ws3.Range("I2").Value = "$newFile1.md"
$ws3.Range("J2").Value = "$newFile1.075c9c4432c64c1ca796c3749c5def87881cb935.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-23 19:05:46"
$ws3.Range("L2").Value = ""
$ws3.Range("M2").Value = "True"
$ws3.Range("N2").Value = ""
$ws3.Range("O2").Value = "False"
$ws3.Range("P2").Value = ""

$ws3.Range("A3").Value = "$newFile2.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "$newFile2.419863fc8d9492e3b6dc1f1704fb5fb5fff97a10.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-23 19:05:26"
$ws3.Range("I3").Value = "$newFile2.md"
$ws3.Range("J3").Value = "$newFile2.419863fc8d9492e3b6dc1f1704fb5fb5fff97a10.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-23 19:05:46"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks: rebuild in final order (A2, I2, A3, I3)
$ws3.Range("A2").Hyperlinks.Delete() | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c916a8f75a7e29c1e11910abb3d50b39a646534e/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c916a8f75a7e29c1e11910abb3d50b39a646534e/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md") | Out-Null

Write-Output "done"
